$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Swap the URL values between B2 and B3 (point back to the Nimbus instance of PPM)
$b2 = $ws.Range("B2").Value()
$b3 = $ws.Range("B3").Value()
$ws.Range("B2").Value = $b3
$ws.Range("B3").Value = $b2
